# Apply updated "dSF" (column F) values as part of a data repull / mean
# recalculation pass. Only column F values change; all other columns and
# the underlying E ("dS0") values stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 5
    7  = 1
    9  = 1
    11 = 0
    14 = -2
    16 = -1
    25 = 3
    33 = 1
    34 = -1
    46 = -3
    47 = -3
    49 = 1
    53 = 3
    54 = -2
    56 = -1
    60 = 2
    68 = -6
    69 = 3
    71 = -11
    73 = -2
    74 = -2
    79 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
